# Update "想去人数" (want-to-go count) figures in column F across the four
# sheets of the 上海-漫展信息 workbook, matching the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 161
$ws.Range("F3").Value = 2419
$ws.Range("F5").Value = 141
$ws.Range("F6").Value = 68
$ws.Range("F7").Value = 279
$ws.Range("F8").Value = 339
$ws.Range("F9").Value = 3262
$ws.Range("F10").Value = 1162
$ws.Range("F15").Value = 1488
$ws.Range("F16").Value = 754
$ws.Range("F17").Value = 1708
$ws.Range("F20").Value = 70
$ws.Range("F23").Value = 2632

# 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 37
$ws.Range("F28").Value = 184
$ws.Range("F39").Value = 245

# 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 2502
$ws.Range("F8").Value = 137
$ws.Range("F9").Value = 7
$ws.Range("F11").Value = 360
$ws.Range("F12").Value = 2795
$ws.Range("F14").Value = 677

# 全部类型 (All types - combined view)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 137
$ws.Range("F4").Value = 7
$ws.Range("F5").Value = 161
$ws.Range("F6").Value = 2419
$ws.Range("F7").Value = 360
$ws.Range("F10").Value = 37
$ws.Range("F13").Value = 68
$ws.Range("F14").Value = 279
$ws.Range("F15").Value = 339
$ws.Range("F16").Value = 1162
$ws.Range("F23").Value = 1488
$ws.Range("F27").Value = 754
$ws.Range("F31").Value = 1708
$ws.Range("F37").Value = 184
$ws.Range("F39").Value = 70
$ws.Range("F45").Value = 2632
